# Insert a new row at row 813, shifting existing rows 813:854 down to 814:855,
# and populate the new row with the 2026/02/13 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(813).Insert()

$ws.Cells.Item(813, 1).Value = "'2026/02/13"
$ws.Cells.Item(813, 2).Value = "金"
$ws.Cells.Item(813, 3).Value = 10
$ws.Cells.Item(813, 4).Value = 201
